# Add a "Todo list" / Bibliography heading before the {@bibliography}
# merge field, and give the merge field's own paragraph a dedicated
# "References" paragraph style (based on Normal, 9pt).

$d = $word.ActiveDocument

# 1. Define the new "References" paragraph style.
$refStyle = $d.Styles.Add("References", 1)
$refStyle.BaseStyle = $d.Styles("Normal")
$refStyle.QuickStyle = $true
$refStyle.Font.Size = 9

# 2. Locate the paragraph containing the {@bibliography} merge field -
#    it is the last paragraph in the document body.
$count = $d.Paragraphs.Count
$biblioIndex = $count
$biblioPara = $d.Paragraphs($biblioIndex)

# 3. Insert a new, empty paragraph immediately before it; this will
#    become the "Bibliography" heading paragraph. After the insert the
#    new blank paragraph takes the old index and the bibliography
#    paragraph's index shifts down by one.
$biblioPara.Range.InsertParagraphBefore()
$headingIndex = $biblioIndex
$biblioIndex = $biblioIndex + 1

# 4. Fill in the heading text/style.
$headingPara = $d.Paragraphs($headingIndex)
$headingPara.Range.Text = "Bibliography"
$headingPara.Style = "Heading2"

# 5. Apply the new "References" style to the {@bibliography} paragraph.
$biblioPara = $d.Paragraphs($biblioIndex)
$biblioPara.Style = "References"
